$wb = $excel.ActiveWorkbook

# --- Sheet 1 (AssessmentResult): insert a new header row above row 1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(1).Insert()

# New header labels land on the row that was just inserted (row 1).
$ws1.Range("B1").Value = "0ne"
$ws1.Range("C1").Value = "tw0"
$ws1.Range("D1").Value = "thr3e"
$ws1.Range("E1").Value = "f0ur"
$ws1.Range("F1").Value = "fiv3"
$ws1.Range("G1").Value = "s1x"
$ws1.Range("H1").Value = "s3ven"

# Restore the previously-selected cell on sheet 1.
$ws1.Range("C4").Select()

# --- Sheet 2 (Criteria): update selection + page setup ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A21").Select()
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
